$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rutas_registros")

$ws.Range("A4").Value = "20240713"
$ws.Range("B4").Value = "ruta ejemplo"

$ws.Range("A5").Value = "20240714"
$ws.Range("B5").Value = "ruta ejemplo"

$ws.Range("A6").Value = "20240702"
$ws.Range("B6").Value = "ruta ejemplo"
